$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.678.14"
$ws.Range("E2").Value = "  +0.83%  "
$ws.Range("D3").Value = "3.503.83"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'603.98"
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").Value = "'170.89"
$ws.Range("E6").Value = "  -1.90%  "
$ws.Range("E7").Value = "  -1.58%  "
$ws.Range("D8").Value = "3.501.89"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'0.199"
$ws.Range("E10").Value = "  +5.14%  "
$ws.Range("E11").Value = "  -0.40%  "
$ws.Range("E12").Value = "  -2.71%  "
$ws.Range("D13").Value = "'47.04"
$ws.Range("E13").Value = "  -0.09%  "
$ws.Range("E14").Value = "  +0.73%  "
$ws.Range("D15").Value = "4.075.94"
$ws.Range("E15").Value = "  +0.26%  "
$ws.Range("D16").Value = "'619.61"
$ws.Range("E16").Value = "  -8.28%  "
$ws.Range("D17").Value = "'8.38"
$ws.Range("E17").Value = "  -3.83%  "
$ws.Range("D18").Value = "3.505.55"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("D19").Value = "69.739.28"
$ws.Range("E19").Value = "  +1.03%  "
$ws.Range("E20").Value = "  -2.11%  "
$ws.Range("E21").Value = "  -1.11%  "
$ws.Range("D22").Value = "'0.881"
$ws.Range("E22").Value = "  -2.43%  "
$ws.Range("D23").Value = "'9.87"
$ws.Range("E23").Value = "  -11.72%  "
$ws.Range("D24").Value = "'15.75"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").Value = "'95.90"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("E28").Value = "  -2.65%  "
$ws.Range("D29").Value = "'9.15"
$ws.Range("E29").Value = "  -3.01%  "
$ws.Range("D30").Value = "'33.04"
$ws.Range("E30").Value = "  +0.53%  "
$ws.Range("E31").Value = "  -4.04%  "
$ws.Range("E32").Value = "  -4.91%  "
$ws.Range("D33").Value = "'1.33"
$ws.Range("E33").Value = "  -1.34%  "
$ws.Range("D34").Value = "'6.93"
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").Value = "'563.86"
$ws.Range("E35").Value = "  -4.52%  "
$ws.Range("D36").Value = "'10.72"
$ws.Range("E36").Value = "  -1.56%  "
$ws.Range("E37").Value = "  -3.57%  "
$ws.Range("D38").Value = "'57.01"
$ws.Range("E38").Value = "  -0.41%  "
$ws.Range("E39").Value = "  -3.75%  "
$ws.Range("E40").Value = "  +0.18%  "
$ws.Range("E41").Value = "  +3.79%  "
$ws.Range("D42").Value = "'0.0446"
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("D43").Value = "3.323.56"
$ws.Range("E43").Value = "  -2.95%  "
$ws.Range("E44").Value = "  -3.60%  "
$ws.Range("B45").Value = "PEPE"
$ws.Range("C45").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D45").Value = "0.0₃0706"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("B46").Value = "ThetaToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D46").Value = "'2.97"
$ws.Range("E46").Value = "  +1.91%  "
$ws.Range("B47").Value = "InjectiveProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D47").Value = "'32.93"
$ws.Range("E47").Value = "  -1.39%  "
$ws.Range("E48").Value = "  +0.28%  "
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").Value = "'135.20"
$ws.Range("E50").Value = "  +2.40%  "
$ws.Range("E51").Value = "  -1.14%  "
